$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 2; $i -le 6; $i++) {
    $num = "{0:D2}" -f ($i - 2)
    $ws.Range("A$i").Value = "sequences/278857_motortesting_sequence_$num.csv"
}
